$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New invoice numbers for column A, rows 2 through 17
$values = @(839974, 860620, 982436, 694088, 305290, 214298, 710472, 427043, 725704, 333552, 701177, 948269, 810093, 706430, 101584, 407478)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# The source sheet previously had a stale selection at C11 (outside the
# new data range); reset it to A1 so the saved view no longer references
# a cell beyond the refreshed data.
$ws.Range("A1").Select()
